$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "2016-08-18 12:14:48" -> "2016-08-18 12:15:42" (shared string also used by
# de-de's Correspond Handoff Datetime column, which happens to hold the same
# timestamp for these two rows)
$wsOverview.Range("G2").Value = "2016-08-18 12:15:42"
$wsOverview.Range("G5").Value = "2016-08-18 12:15:42"
$wsDeDe.Range("H2").Value = "2016-08-18 12:15:42"
$wsDeDe.Range("H5").Value = "2016-08-18 12:15:42"

# Priority column (E) ht -> mt (shared across zh-cn and de-de sheets)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-18 12:15:37"
$wsZhCn.Range("H5").Value = "2016-08-18 12:15:37"

# zh-cn sheet: Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-18 12:16:09"
$wsZhCn.Range("K5").Value = "2016-08-18 12:16:09"

# de-de sheet: Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-18 12:16:17"
$wsDeDe.Range("K5").Value = "2016-08-18 12:16:17"
